# Auto-generated Excel COM-interop script to apply value updates
# described by the Kraken_Profits.xlsx diff (per-sheet currentAveragePrice
# / Leve profit recalculation columns H:N).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 83.14286
$ws.Range("I2").Value = 85.5
$ws.Range("K2").Value = 85.5
$ws.Range("M2").Value = 27.5

$ws.Range("H18").Value = 2366.6667
$ws.Range("I18").Value = 2366.6667
$ws.Range("K18").Value = 2366.6667
$ws.Range("M18").Value = -2082.6667

$ws.Range("H19").Value = 737
$ws.Range("I19").Value = 923.5
$ws.Range("J19").Value = 550.5
$ws.Range("K19").Value = 923.5
$ws.Range("L19").Value = 550.5
$ws.Range("M19").Value = -748.5
$ws.Range("N19").Value = -900.5

$ws.Range("H32").Value = 8468.1875
$ws.Range("J32").Value = 9444.223
$ws.Range("L32").Value = 9444.223
$ws.Range("N32").Value = -10096.223

$ws.Range("H38").Value = 748.1111
$ws.Range("I38").Value = 33.25
$ws.Range("J38").Value = 1320
$ws.Range("K38").Value = 99.75
$ws.Range("L38").Value = 3960
$ws.Range("M38").Value = 272.25
$ws.Range("N38").Value = -4704

$ws.Range("H40").Value = 5366.5835
$ws.Range("I40").Value = 2199.8333
$ws.Range("J40").Value = 8533.333000000001
$ws.Range("K40").Value = 2199.8333
$ws.Range("L40").Value = 8533.333000000001
$ws.Range("M40").Value = -2024.8333
$ws.Range("N40").Value = -8883.333000000001

$ws.Range("H61").Value = 5016.5
$ws.Range("J61").Value = 5016.5
$ws.Range("L61").Value = 15049.5
$ws.Range("N61").Value = -15393.5

$ws.Range("H64").Value = 4200
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 4200
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 4200
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4696

$ws.Range("H67").Value = 4200
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 4200
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 4200
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -5916

$ws.Range("H86").Value = 7711
$ws.Range("I86").Value = 6185
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 6185
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -5062
$ws.Range("N86").Value = -12246

$ws.Range("H89").Value = 7711
$ws.Range("I89").Value = 6185
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 30925
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -25309
$ws.Range("N89").Value = -61232

$ws.Range("H92").Value = 1252.6923
$ws.Range("I92").Value = 1117.1818
$ws.Range("J92").Value = 1998
$ws.Range("K92").Value = 1117.1818
$ws.Range("L92").Value = 1998
$ws.Range("M92").Value = 130.8181999999999
$ws.Range("N92").Value = -4494

$ws.Range("H98").Value = 397.44446
$ws.Range("I98").Value = 421.33334
$ws.Range("K98").Value = 421.33334
$ws.Range("M98").Value = 1076.66666

$ws.Range("H113").Value = 966.6667
$ws.Range("I113").Value = 966.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 966.6667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 2287.3333
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 397.44446
$ws.Range("I122").Value = 421.33334
$ws.Range("K122").Value = 1264.00002
$ws.Range("M122").Value = 1185.99998

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H137").Value = 3257.3462
$ws.Range("I137").Value = 2461.375
$ws.Range("K137").Value = 7384.125
$ws.Range("M137").Value = -4834.125

$ws.Range("H141").Value = 49999.5
$ws.Range("I141").Value = 49999.5
$ws.Range("K141").Value = 149998.5
$ws.Range("M141").Value = -144818.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3550.6365
$ws.Range("I110").Value = 1211.6
$ws.Range("J110").Value = 5499.8335
$ws.Range("K110").Value = 1211.6
$ws.Range("L110").Value = 5499.8335
$ws.Range("M110").Value = 833.4000000000001
$ws.Range("N110").Value = -9589.833500000001

$ws.Range("H122").Value = 3489.2
$ws.Range("I122").Value = 4223
$ws.Range("K122").Value = 12669
$ws.Range("M122").Value = -10219

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 99780
$ws.Range("J132").Value = 99780
$ws.Range("L132").Value = 99780
$ws.Range("N132").Value = -109900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2995
$ws.Range("I58").Value = 2995
$ws.Range("K58").Value = 2995
$ws.Range("M58").Value = -2792

$ws.Range("H94").Value = 1994.6
$ws.Range("J94").Value = 1996.3334
$ws.Range("L94").Value = 1996.3334
$ws.Range("N94").Value = -2898.3334

$ws.Range("H99").Value = 1986.5
$ws.Range("I99").Value = 1986.5
$ws.Range("K99").Value = 1986.5
$ws.Range("M99").Value = -488.5

$ws.Range("H122").Value = 1582.3334
$ws.Range("I122").Value = 1582.3334
$ws.Range("K122").Value = 4747.0002
$ws.Range("M122").Value = -2297.0002

$ws.Range("H126").Value = 1986.5
$ws.Range("I126").Value = 1986.5
$ws.Range("K126").Value = 5959.5
$ws.Range("M126").Value = -3489.5

$ws.Range("H132").Value = 5221.25
$ws.Range("I132").Value = 5478.6665
$ws.Range("J132").Value = 4963.8335
$ws.Range("K132").Value = 16435.9995
$ws.Range("L132").Value = 14891.5005
$ws.Range("M132").Value = -13905.9995
$ws.Range("N132").Value = -19951.5005

$ws.Range("H136").Value = 2995
$ws.Range("I136").Value = 2995
$ws.Range("K136").Value = 8985
$ws.Range("M136").Value = -6435

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 16.5
$ws.Range("I12").Value = 9
$ws.Range("J12").Value = 18
$ws.Range("K12").Value = 27
$ws.Range("L12").Value = 54
$ws.Range("M12").Value = 146
$ws.Range("N12").Value = -400

$ws.Range("H128").Value = 89000
$ws.Range("I128").Value = 89000
$ws.Range("K128").Value = 267000
$ws.Range("M128").Value = -262020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 9142.857
$ws.Range("I43").Value = 1333.3334
$ws.Range("K43").Value = 1333.3334
$ws.Range("M43").Value = -1182.3334

$ws.Range("H122").Value = 9739.200000000001
$ws.Range("J122").Value = 11332.667
$ws.Range("L122").Value = 33998.001
$ws.Range("N122").Value = -38898.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4300
$ws.Range("I16").Value = 4300
$ws.Range("K16").Value = 4300
$ws.Range("M16").Value = -4130

$ws.Range("H46").Value = 825
$ws.Range("I46").Value = 791.6667
$ws.Range("J46").Value = 875
$ws.Range("K46").Value = 791.6667
$ws.Range("L46").Value = 875
$ws.Range("M46").Value = -603.6667
$ws.Range("N46").Value = -1251

$ws.Range("H101").Value = 15787.333
$ws.Range("J101").Value = 15787.333
$ws.Range("L101").Value = 15787.333
$ws.Range("N101").Value = -22277.333

$ws.Range("H132").Value = 34331
$ws.Range("I132").Value = 39397.4
$ws.Range("K132").Value = 118192.2
$ws.Range("M132").Value = -115662.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 30013
$ws.Range("I22").Value = 30013
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 30013
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -29720
$ws.Range("N22").ClearContents()

$ws.Range("H126").Value = 3801.6
$ws.Range("I126").Value = 3252
$ws.Range("K126").Value = 9756
$ws.Range("M126").Value = -7286

$ws.Range("H132").Value = 5996.143
$ws.Range("I132").Value = 5996.143
$ws.Range("K132").Value = 17988.429
$ws.Range("M132").Value = -15458.429

$ws.Range("H136").Value = 4593.6665
$ws.Range("I136").Value = 4593.6665
$ws.Range("K136").Value = 13780.9995
$ws.Range("M136").Value = -11230.9995
